$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Set B7 ("Experimental" row) to the literal text "true" (not a boolean).
# A direct Value assignment of "true"/"false" is auto-coerced to a Boolean
# by Excel, so instead write it as a text formula and then convert the
# formula to a static value via copy / paste-special (values only) -
# this keeps the cell's existing style and produces a genuine string cell.
$c = $ws.Cells.Item(7, 2)
$c.Formula = "=""true"""
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues

# Update the "Date" value (B8) to the new generation timestamp.
$ws.Cells.Item(8, 2).Value = "2023-02-01T09:05:11-06:00"
